# Deploying to gh-pages: extend the debt-service table from 2022 through 2024
# (adds columns T:U, mirrors the existing row 3/4/5 formatting into them,
# widens the new year columns, and grows row 5 to fit the extra wrapped text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year columns (T=2023, U=2024) need the same width as the rest of the
# year columns (D:U). (COM ColumnWidth snaps to whole-pixel steps, so this is
# the closest achievable value to the authored 8.7109375 character-width.)
$ws.Range("D1:U1").EntireColumn.ColumnWidth = 7.83

# Bring over the formatting (styles/borders) from column S into the two new
# columns for the header block (rows 3-5), then fill in the actual values.
$ws.Range("S3:S5").Copy()
$ws.Range("T3:T5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S3:S5").Copy()
$ws.Range("U3:U5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Year header row
$ws.Cells.Item(4, 20).Value = 2023
$ws.Cells.Item(4, 21).Value = 2024

# Data row
$ws.Cells.Item(5, 20).Value = 10.8
$ws.Cells.Item(5, 21).Value = 6.5

# Row 5 grows a little taller to fit the (now wider) wrapped header text.
$ws.Rows.Item(5).RowHeight = 41.25
